# Generate Report for Handoff
# Updates the handoff-generation timestamps and sets the Priority column
# ("ht") for the rows that correspond to file
# 48525296-5182-4abf-addc-60195e3bf9eb.md, now that a handoff report has
# been generated for it.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Rows 9-14 on each sheet correspond to the same six files; update them all.
for ($row = 9; $row -le 14; $row++) {
    # zh-cn / de-de: Priority column (E) becomes "ht" (handoff type)
    $zhcn.Range("E$row").Value = "ht"
    $dede.Range("E$row").Value = "ht"
}

# Overview sheet: Latest HO Xliff Generate Date (column G) updated
$overview.Range("G9:G14").Value = "2016-09-04 02:24:07"

# zh-cn sheet: Latest Handoff Datetime (column H) updated
$zhcn.Range("H9:H14").Value = "2016-09-04 02:23:58"
